$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "BORGASONE TRIO  CREAM" (row 6) - entire row shifts up
$ws.Rows.Item(6).Delete()

# Renumber the "م" (sequence number) column for the rows that shifted up (rows 6-18)
for ($r = 6; $r -le 18; $r++) {
    $ws.Range("A$r").Value = $r - 3
}

# Update the total (سعر البيع sum) to reflect the removed row
$ws.Range("K19").Value = 344.5
